# SSDM-13163: Added the description attribute to exports.
# This script rewrites the "type" metadata block (rows 2-3) to insert a new
# "Description" column, and converts the boolean Mandatory/"Show in edit
# views" columns (and the type-level "Auto generate codes" flag) from real
# booleans into plain text "TRUE"/"FALSE" cells (background-colour based
# assertions on booleans were dropped from the tests, so the exporter now
# emits these as text).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$boolFormat = '"TRUE";"TRUE";"FALSE"'

# ---------------------------------------------------------------------
# Row 2 (sample type header labels) + Row 3 (sample type values):
# insert a new "Description" column at C, shifting the former C/D/E
# content to D/E/F.
# ---------------------------------------------------------------------

# Row 2 labels
$ws.Range("F2").Value = "Generated code prefix"
$ws.Range("E2").Value = "Validation script"
$ws.Range("D2").Value = "Auto generate codes"
$ws.Range("C2").Value = "Description"

$ws.Range("C2").Font.Name = "Calibri"
$ws.Range("C2").Font.Size = 14
$ws.Range("C2").Font.Bold = $true
$ws.Range("C2").HorizontalAlignment = -4131

$ws.Range("D2").Font.Name = "Calibri"
$ws.Range("D2").Font.Size = 14
$ws.Range("D2").Font.Bold = $true
$ws.Range("D2").HorizontalAlignment = -4131

$ws.Range("E2").Font.Name = "Calibri"
$ws.Range("E2").Font.Size = 14
$ws.Range("E2").Font.Bold = $true
$ws.Range("E2").HorizontalAlignment = 1

$ws.Range("F2").Font.Name = "Calibri"
$ws.Range("F2").Font.Size = 13
$ws.Range("F2").Font.Bold = $true
$ws.Range("F2").HorizontalAlignment = 1

# Row 3 values
$ws.Range("E3").Value = "date_range_validation.py"
$ws.Range("D3").NumberFormat = $boolFormat
$ws.Range("D3").Value = "'FALSE"
$ws.Range("C3").Value = "Course"

$ws.Range("C3").Font.Name = "Calibri"
$ws.Range("C3").Font.Size = 12
$ws.Range("C3").HorizontalAlignment = 1

$ws.Range("D3").Font.Name = "Calibri"
$ws.Range("D3").Font.Size = 12
$ws.Range("D3").HorizontalAlignment = -4131

$ws.Range("E3").Font.Name = "Calibri"
$ws.Range("E3").Font.Size = 11
$ws.Range("E3").HorizontalAlignment = 1

# ---------------------------------------------------------------------
# Property-assignment rows (5-9): Mandatory / "Show in edit views" move
# from real booleans to text "TRUE"/"FALSE" (same displayed value, but
# stored as a string so it survives the cross-OS comparison).
# ---------------------------------------------------------------------

foreach ($r in 5..9) {
    foreach ($col in @("C", "D")) {
        $cell = $ws.Range("$col$r")
        $cell.NumberFormat = $boolFormat
    }
}

$ws.Range("C5").Value = "'TRUE"
$ws.Range("D5").Value = "'TRUE"

$ws.Range("C6").Value = "'TRUE"
$ws.Range("D6").Value = "'TRUE"

$ws.Range("C7").Value = "'TRUE"
$ws.Range("D7").Value = "'TRUE"

$ws.Range("C8").Value = "'TRUE"
$ws.Range("D8").Value = "'TRUE"

$ws.Range("C9").Value = "'FALSE"
$ws.Range("D9").Value = "'TRUE"

# ---------------------------------------------------------------------
# Row heights tweak slightly after the re-export
# ---------------------------------------------------------------------
$ws.Rows.Item(2).RowHeight = 17.35
$ws.Rows.Item(3).RowHeight = 15

# ---------------------------------------------------------------------
# Two trailing blank rows now appear after the re-export (13 & 14)
# ---------------------------------------------------------------------
$ws.Cells.Item(13, 1).RowHeight = 15
$ws.Cells.Item(14, 1).RowHeight = 15

# ---------------------------------------------------------------------
# Selection moves to D2:D3 with D2 active
# ---------------------------------------------------------------------
$ws.Range("D2:D3").Select()
